$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad" date) for rows 2-23 from 45207 to 45208
for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}
